$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 491.05554
$ws.Range("J17").Value = 491.05554
$ws.Range("L17").Value = 1473.16662
$ws.Range("N17").Value = -1809.16662

$ws.Range("H38").Value = 244.42857
$ws.Range("I38").Value = 244.42857
$ws.Range("K38").Value = 733.28571
$ws.Range("M38").Value = -361.28571

$ws.Range("H43").Value = 4281804.5
$ws.Range("I43").Value = 34467
$ws.Range("J43").Value = 5556006
$ws.Range("K43").Value = 34467
$ws.Range("L43").Value = 5556006
$ws.Range("M43").Value = -34398
$ws.Range("N43").Value = -5556144

$ws.Range("H74").Value = 3224.125
$ws.Range("I74").Value = 2965.5
$ws.Range("K74").Value = 2965.5
$ws.Range("M74").Value = -2029.5

$ws.Range("H77").Value = 3224.125
$ws.Range("I77").Value = 2965.5
$ws.Range("K77").Value = 14827.5
$ws.Range("M77").Value = -10147.5

$ws.Range("H92").Value = 679.8421
$ws.Range("I92").Value = 553.94116
$ws.Range("K92").Value = 553.94116
$ws.Range("M92").Value = 694.05884

$ws.Range("H112").Value = 2495.5144
$ws.Range("J112").Value = 3122.8462
$ws.Range("L112").Value = 9368.5386
$ws.Range("N112").Value = -11584.5386

$ws.Range("H132").Value = 10106686
$ws.Range("J132").Value = 2933.3333
$ws.Range("L132").Value = 8799.999899999999
$ws.Range("N132").Value = -13859.9999

$ws.Range("H135").Value = 47619988
$ws.Range("I135").Value = 453.7857
$ws.Range("K135").Value = 4084.0713
$ws.Range("M135").Value = -1549.0713

$ws.Range("H138").Value = 1421.38
$ws.Range("I138").Value = 616.6667
$ws.Range("J138").Value = 1935.8689
$ws.Range("K138").Value = 1850.0001
$ws.Range("L138").Value = 5807.6067
$ws.Range("M138").Value = 3289.9999
$ws.Range("N138").Value = -16087.6067

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3313.2134
$ws.Range("I32").Value = 2987.1384
$ws.Range("J32").Value = 5432.7
$ws.Range("K32").Value = 2987.1384
$ws.Range("L32").Value = 5432.7
$ws.Range("M32").Value = -2700.1384
$ws.Range("N32").Value = -6006.7

$ws.Range("H35").Value = 2180.3333
$ws.Range("I35").Value = 2180.3333
$ws.Range("K35").Value = 2180.3333
$ws.Range("M35").Value = -1774.3333

$ws.Range("H45").Value = 1056.08
$ws.Range("I45").Value = 1023.64703
$ws.Range("J45").Value = 1125
$ws.Range("K45").Value = 1023.64703
$ws.Range("L45").Value = 1125
$ws.Range("M45").Value = -646.64703
$ws.Range("N45").Value = -1879

$ws.Range("H61").Value = 1137.4634
$ws.Range("I61").Value = 992.51514
$ws.Range("J61").Value = 1735.375
$ws.Range("K61").Value = 992.51514
$ws.Range("L61").Value = 1735.375
$ws.Range("M61").Value = -780.51514
$ws.Range("N61").Value = -2159.375

$ws.Range("H122").Value = 915.85
$ws.Range("I122").Value = 935.1177
$ws.Range("J122").Value = 806.6667
$ws.Range("K122").Value = 2805.3531
$ws.Range("L122").Value = 2420.0001
$ws.Range("M122").Value = -355.3531000000003
$ws.Range("N122").Value = -7320.0001

$ws.Range("H136").Value = 1137.4634
$ws.Range("I136").Value = 992.51514
$ws.Range("J136").Value = 1735.375
$ws.Range("K136").Value = 2977.54542
$ws.Range("L136").Value = 5206.125
$ws.Range("M136").Value = -427.5454199999999
$ws.Range("N136").Value = -10306.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3970.4146
$ws.Range("I134").Value = 1033.1428
$ws.Range("J134").Value = 10296.846
$ws.Range("K134").Value = 3099.4284
$ws.Range("L134").Value = 30890.538
$ws.Range("M134").Value = -564.4284000000002
$ws.Range("N134").Value = -35960.538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1110.7966
$ws.Range("I31").Value = 1126.2808
$ws.Range("J31").Value = 669.5
$ws.Range("K31").Value = 1126.2808
$ws.Range("L31").Value = 669.5
$ws.Range("M31").Value = -831.2808
$ws.Range("N31").Value = -1259.5

$ws.Range("H34").Value = 1110.7966
$ws.Range("I34").Value = 1126.2808
$ws.Range("J34").Value = 669.5
$ws.Range("K34").Value = 1126.2808
$ws.Range("L34").Value = 669.5
$ws.Range("M34").Value = -924.2808
$ws.Range("N34").Value = -1073.5

$ws.Range("H134").Value = 1139.6842
$ws.Range("I134").Value = 1015.8
$ws.Range("K134").Value = 3047.4
$ws.Range("M134").Value = -512.3999999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 6990.909
$ws.Range("I96").Value = 5000
$ws.Range("J96").Value = 7190
$ws.Range("K96").Value = 15000
$ws.Range("L96").Value = 21570
$ws.Range("M96").Value = -12941
$ws.Range("N96").Value = -25688

$ws.Range("H136").Value = 1441.6842
$ws.Range("I136").Value = 917.3333
$ws.Range("K136").Value = 2751.9999
$ws.Range("M136").Value = 2348.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 64288500
$ws.Range("I70").Value = 62502924
$ws.Range("K70").Value = 62502924
$ws.Range("M70").Value = -62502654

$ws.Range("H73").Value = 64288500
$ws.Range("I73").Value = 62502924
$ws.Range("K73").Value = 62502924
$ws.Range("M73").Value = -62501988

$ws.Range("H132").Value = 2011.425
$ws.Range("I132").Value = 1430.9166
$ws.Range("K132").Value = 4292.7498
$ws.Range("M132").Value = -1762.7498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1963.5714
$ws.Range("I7").Value = 1755
$ws.Range("J7").Value = 2339
$ws.Range("K7").Value = 1755
$ws.Range("L7").Value = 2339
$ws.Range("M7").Value = -1643
$ws.Range("N7").Value = -2563

$ws.Range("H22").Value = 2108.1667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2108.1667
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = 2108.1667
$ws.Range("N22").Value = -2698.1667
$ws.Range("L22").ClearContents()

$ws.Range("H27").Value = 2108.1667
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2108.1667
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = 2108.1667
$ws.Range("N27").Value = -2322.1667
$ws.Range("L27").ClearContents()

$ws.Range("H35").Value = 1533
$ws.Range("J35").Value = 800
$ws.Range("L35").Value = 800
$ws.Range("N35").Value = -1472

$ws.Range("H40").Value = 2802.5334
$ws.Range("I40").Value = 1784.3
$ws.Range("K40").Value = 1784.3
$ws.Range("M40").Value = -1648.3

$ws.Range("H43").Value = 933.3333
$ws.Range("J43").Value = 800
$ws.Range("L43").Value = 800
$ws.Range("N43").Value = -1186

$ws.Range("H69").Value = 35833.332
$ws.Range("I69").Value = 25000
$ws.Range("J69").Value = 41250
$ws.Range("K69").Value = 25000
$ws.Range("L69").Value = 41250
$ws.Range("M69").Value = -24189
$ws.Range("N69").Value = -42872

$ws.Range("H72").Value = 35833.332
$ws.Range("I72").Value = 25000
$ws.Range("J72").Value = 41250
$ws.Range("K72").Value = 75000
$ws.Range("L72").Value = 123750
$ws.Range("M72").Value = -70944
$ws.Range("N72").Value = -131862

$ws.Range("H126").Value = 1963.5714
$ws.Range("I126").Value = 1755
$ws.Range("J126").Value = 2339
$ws.Range("K126").Value = 5265
$ws.Range("L126").Value = 7017
$ws.Range("M126").Value = -2795
$ws.Range("N126").Value = -11957

$ws.Range("H136").Value = 1988.2222
$ws.Range("I136").Value = 2761.3333
$ws.Range("J136").Value = 1601.6666
$ws.Range("K136").Value = 8283.999899999999
$ws.Range("L136").Value = 4804.9998
$ws.Range("M136").Value = -5733.999899999999
$ws.Range("N136").Value = -9904.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2226

$ws.Range("H81").Value = 347.14285
$ws.Range("I81").Value = 286
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 572
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = 489
$ws.Range("N81").Value = -3122

$ws.Range("H84").Value = 347.14285
$ws.Range("I84").Value = 286
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 2860
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = 2444
$ws.Range("N84").Value = -15608

$ws.Range("H132").Value = 1151.9318
$ws.Range("I132").Value = 1024.9062
$ws.Range("J132").Value = 1490.6666
$ws.Range("K132").Value = 3074.7186
$ws.Range("L132").Value = 4471.9998
$ws.Range("M132").Value = -544.7185999999997
$ws.Range("N132").Value = -9531.9998

$ws.Range("H136").Value = 1374.6666
$ws.Range("I136").Value = 704
$ws.Range("J136").Value = 1508.8
$ws.Range("K136").Value = 2112
$ws.Range("L136").Value = 4526.4
$ws.Range("M136").Value = 438
$ws.Range("N136").Value = -9626.4
